$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Đơn sale chính" ---
# Delete columns (rightmost-first so shifting doesn't disturb the ones still to delete):
#   AA (Chiết khấu sale phụ), Y (Tỉ lệ chiết khấu sale phụ), W (Công phụ phẫu 2),
#   V (Công phụ phẫu 1), U (Phụ phẫu 2), T (Phụ phẫu 1), S (Bác sĩ 2), R (Bác sĩ 1),
#   Q (Dư nợ), O (Trả sau), N (Thanh toán lần đầu), I (Sale chính), G (Nhóm dịch vụ)
$ws1 = $wb.Worksheets.Item("Đơn sale chính")
$cols1 = @("AA","Y","W","V","U","T","S","R","Q","O","N","I","G")
foreach ($c in $cols1) {
    $ws1.Range($c + "1").EntireColumn.Delete()
}

# --- Sheet 2: "Đơn phụ phẫu 1" ---
# Delete columns (rightmost-first):
#   AA, Z, Y, X, W, U, S, R, Q, P, O, N, M, L, K, J, I, G
$ws2 = $wb.Worksheets.Item("Đơn phụ phẫu 1")
$cols2 = @("AA","Z","Y","X","W","U","S","R","Q","P","O","N","M","L","K","J","I","G")
foreach ($c in $cols2) {
    $ws2.Range($c + "1").EntireColumn.Delete()
}

# --- Sheet 3: "Lương" ---
$ws3 = $wb.Worksheets.Item("Lương")
$ws3.Range("A1").Value = "Danh mục lương"
$ws3.Range("B2").Value = 18
$ws3.Range("B3").Value = 630000
$ws3.Range("B12").Value = 1928571.428571429
$ws3.Range("B29").Value = 2758571.428571429
$ws3.Range("B31").Value = 2758571.428571429
